$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(74).Insert()

$ws.Cells.Item(74,1).Value = 9
$ws.Cells.Item(74,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(74,3).Value = "Metropolitana"
$ws.Cells.Item(74,4).Value = 45280
$ws.Cells.Item(74,5).Value = 13
$ws.Cells.Item(74,6).Value = 100112005
$ws.Cells.Item(74,7).Value = "Puerro"
$ws.Cells.Item(74,8).Value = "Sin especificar"
$ws.Cells.Item(74,9).Value = "Primera"
$ws.Cells.Item(74,10).Value = 70
$ws.Cells.Item(74,11).Value = 7000
$ws.Cells.Item(74,12).Value = 8000
$ws.Cells.Item(74,13).Value = 7500
$ws.Cells.Item(74,14).Value = "$/paquete 20 unidades"
$ws.Cells.Item(74,15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(74,16).Value = 375
$ws.Cells.Item(74,17).Value = 20
$ws.Cells.Item(74,18).Value = "Hortaliza"
